$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update parameter values in row 2
$ws.Range("C2").Value = 100
$ws.Range("E2").Value = 25
$ws.Range("F2").Value = 60000

# Recalculate formula-dependent cell(s)
$wb.Application.Calculate()

# Update the active selection to H3
$ws.Activate()
$ws.Range("H3").Select()
